$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Games")

# Copy the formatting (borders/style) from the row above (C103) down into
# the new rows so the new cells pick up the same "s=10" style Excel would
# normally apply when the row is filled in like its neighbours.
$ws.Range("C103").Copy()
$ws.Range("C104:C106").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Add new game results for 07/04/2019 (FB beat GS, 3 matches) in rows 104-106
$ws.Range("A104").Value = "FB"
$ws.Range("B104").Value = "GS"
$ws.Range("C104").Value = -1

$ws.Range("A105").Value = "FB"
$ws.Range("B105").Value = "GS"
$ws.Range("C105").Value = -1

$ws.Range("A106").Value = "FB"
$ws.Range("B106").Value = "GS"
$ws.Range("C106").Value = -1

$ws.Range("F106").Select()

$wb.Save()
